$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells so numeric-looking strings keep their exact literal text
$ws.Range("D2,D3,D4,D5,D6,D7,D8,D10,D11,D12,D14,D15,D16,D17,D18,D19,D20,D23,D24,D25,D26,D28,D31,D34,D35,D36,D39,D41,D42,D43,D44,D46,D47,D49,D50").NumberFormat = "@"

# Apply updated cell values
$ws.Range("D2").Value = "51.399.97"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").Value = "2.754.67"
$ws.Range("E3").Value = "  +4.63%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "115.85"
$ws.Range("E5").Value = "  +3.27%  "
$ws.Range("D6").Value = "332.48"
$ws.Range("D7").Value = "0.538"
$ws.Range("E7").Value = "  +2.25%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +5.68%  "
$ws.Range("D10").Value = "41.58"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("D11").Value = "0.0857"
$ws.Range("E11").Value = "  +5.73%  "
$ws.Range("D12").Value = "20.24"
$ws.Range("E12").Value = "  +2.41%  "
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "7.62"
$ws.Range("E14").Value = "  +4.89%  "
$ws.Range("D15").Value = "3.184.47"
$ws.Range("E15").Value = "  +4.79%  "
$ws.Range("D16").Value = "2.767.73"
$ws.Range("E16").Value = "  +5.27%  "
$ws.Range("D17").Value = "0.885"
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "51.417.91"
$ws.Range("E18").Value = "  +4.59%  "
$ws.Range("D19").Value = "3.22"
$ws.Range("E19").Value = "  +6.20%  "
$ws.Range("D20").Value = "13.44"
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").Value = "278.05"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("D24").Value = "69.65"
$ws.Range("D25").Value = "2.65"
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("D26").Value = "26.81"
$ws.Range("E26").Value = "  +2.36%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "10.19"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").Value = "35.22"
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "0.0824"
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").Value = "19.07"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("D39").Value = "3.24"
$ws.Range("E39").Value = "  +3.43%  "
$ws.Range("E40").Value = "  +10.55%  "
$ws.Range("D41").Value = "126.72"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "23.07"
$ws.Range("E42").Value = "  +4.36%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "2.30"
$ws.Range("E43").Value = "  +7.98%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "0.114"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("E45").Value = "  +13.19%  "
$ws.Range("D46").Value = "2.088.74"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +3.34%  "
$ws.Range("D49").Value = "5.53"
$ws.Range("E49").Value = "  +6.11%  "
$ws.Range("D50").Value = "9.00"
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("E51").Value = "  +1.77%  "
